$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Rows.Item(6).AutoFit()

$ws.Range("A2:E5").Select()
